# Refresh the crypto price ticker (Price + Volume(1h)) to the latest snapshot.
# Values are kept as literal text (prices like "27.142.61" / "41.90" are not
# valid Excel numbers/would lose formatting), matching the original sheet,
# where column D/E are plain text cells, not numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "27.142.61"; E = "  -2.52%  " }
    @{ Row = 3; D = "1.711.02"; E = "  -2.90%  " }
    @{ Row = 4; D = "0.9997"; E = "  -0.33%  " }
    @{ Row = 5; D = "307.31"; E = "  -6.28%  " }
    @{ Row = 6; D = "0.9996"; E = "  -0.19%  " }
    @{ Row = 7; D = "0.4791"; E = "  +7.37%  " }
    @{ Row = 8; D = "0.3438"; E = "  -2.88%  " }
    @{ Row = 9; D = "41.90"; E = "  -0.09%  " }
    @{ Row = 10; D = "0.07296"; E = "  -1.46%  " }
    @{ Row = 11; D = "1.049"; E = "  -4.45%  " }
    @{ Row = 12; D = "0.9991"; E = "  -0.30%  " }
    @{ Row = 13; D = "19.87"; E = "  -4.76%  " }
    @{ Row = 14; D = "5.855"; E = "  -2.68%  " }
    @{ Row = 15; D = "1.713.30"; E = "  -2.87%  " }
    @{ Row = 16; D = "6.837"; E = "  -5.45%  " }
    @{ Row = 17; D = "89.16"; E = "  -3.92%  " }
    @{ Row = 18; D = "0.00001041"; E = "  -1.74%  " }
    @{ Row = 19; D = "NONE"; E = "  -1.20%  " }
    @{ Row = 20; D = "0.9990"; E = "  -0.21%  " }
    @{ Row = 21; D = "16.45"; E = "  -3.63%  " }
    @{ Row = 22; D = "5.601"; E = "  -2.64%  " }
    @{ Row = 23; D = "27.163.44"; E = "  -2.62%  " }
    @{ Row = 24; D = "NONE"; E = "  -3.47%  " }
    @{ Row = 25; D = "2.092"; E = "  -0.78%  " }
    @{ Row = 26; D = "155.10"; E = "  -3.58%  " }
    @{ Row = 27; D = "19.63"; E = "  -3.57%  " }
    @{ Row = 28; D = "1.904.14"; E = "  -3.20%  " }
    @{ Row = 29; D = "2.083"; E = "  -2.76%  " }
    @{ Row = 30; D = "119.38"; E = "  -3.91%  " }
    @{ Row = 31; D = "1.009"; E = "  -8.09%  " }
    @{ Row = 32; D = "0.09213"; E = "  +0.21%  " }
    @{ Row = 33; D = "3.578"; E = "  -3.02%  " }
    @{ Row = 34; D = "5.311"; E = "  -5.88%  " }
    @{ Row = 35; D = "0.02196"; E = "  -3.60%  " }
    @{ Row = 36; D = "0.05835"; E = "  -5.63%  " }
    @{ Row = 37; D = "NONE"; E = "  -6.42%  " }
    @{ Row = 38; D = "NONE"; E = "  -4.99%  " }
    @{ Row = 39; D = "4.740"; E = "  -4.08%  " }
    @{ Row = 40; D = "0.9995"; E = "  -0.05%  " }
    @{ Row = 41; D = "1.401"; E = "  +0.39%  " }
    @{ Row = 42; D = "0.5875"; E = "  -6.79%  " }
    @{ Row = 43; D = "NONE"; E = "  -6.28%  " }
    @{ Row = 44; D = "7.469"; E = "  -4.96%  " }
    @{ Row = 45; D = "12.62"; E = "  -5.04%  " }
    @{ Row = 46; D = "3.555"; E = "  -5.01%  " }
    @{ Row = 47; D = "0.5625"; E = "  -3.82%  " }
    @{ Row = 48; D = "117.59"; E = "  -3.76%  " }
    @{ Row = 49; D = "1.842"; E = "  -5.53%  " }
    @{ Row = 50; D = "0.06619"; E = "  -3.87%  " }
    @{ Row = 51; D = "1.085"; E = "  -4.23%  " }
)

foreach ($u in $updates) {
    if ($u.D -ne "NONE") {
        $d = $u.D
        # Force literal-text entry for plain-decimal-looking prices (e.g. "0.9997",
        # "41.90") so Excel stores/keeps them as text instead of coercing to a
        # number (which would also strip meaningful trailing/leading zeros).
        # Multi-dot values (e.g. "27.142.61") are never valid numbers so they
        # already round-trip as text without the quote prefix.
        if ($d -match "^[+-]?\d+(\.\d+)?([eE][+-]?\d+)?$") {
            $d = "'" + $d
        }
        $ws.Range("D" + $u.Row).Value = $d
    }
    $ws.Range("E" + $u.Row).Value = $u.E
}
